$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17,8).Value = 2579.5217
$ws.Cells.Item(17,9).Value = 70
$ws.Cells.Item(17,10).Value = 2693.5908
$ws.Cells.Item(17,11).Value = 210
$ws.Cells.Item(17,12).Value = 8080.7724
$ws.Cells.Item(17,13).Value = -42
$ws.Cells.Item(17,14).Value = -8416.7724

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40,8).Value = 2518
$ws.Cells.Item(40,9).Value = 1083.3334
$ws.Cells.Item(40,10).Value = 3235.3333
$ws.Cells.Item(40,11).Value = 1083.3334
$ws.Cells.Item(40,12).Value = 3235.3333
$ws.Cells.Item(40,13).Value = -908.3334
$ws.Cells.Item(40,14).Value = -3585.3333

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113,8).Value = 4169.3335
$ws.Cells.Item(113,9).Value = 3457.5
$ws.Cells.Item(113,10).Value = 4982.857
$ws.Cells.Item(113,11).Value = 3457.5
$ws.Cells.Item(113,12).Value = 4982.857
$ws.Cells.Item(113,13).Value = -203.5
$ws.Cells.Item(113,14).Value = -11490.857

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125,8).Value = 3664.2222
$ws.Cells.Item(125,9).Value = 1329.6666
$ws.Cells.Item(125,10).Value = 4831.5
$ws.Cells.Item(125,11).Value = 11966.9994
$ws.Cells.Item(125,12).Value = 43483.5
$ws.Cells.Item(125,13).Value = -9506.999400000001
$ws.Cells.Item(125,14).Value = -48403.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132,8).Value = 9125.166999999999
$ws.Cells.Item(132,9).Value = 7270.6816
$ws.Cells.Item(132,10).Value = 14225
$ws.Cells.Item(132,11).Value = 21812.0448
$ws.Cells.Item(132,12).Value = 42675
$ws.Cells.Item(132,13).Value = -19282.0448
$ws.Cells.Item(132,14).Value = -47735

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137,8).Value = 2298.054
$ws.Cells.Item(137,9).Value = 2724.7368
$ws.Cells.Item(137,10).Value = 1847.6666
$ws.Cells.Item(137,11).Value = 8174.2104
$ws.Cells.Item(137,12).Value = 5542.9998
$ws.Cells.Item(137,13).Value = -5624.2104
$ws.Cells.Item(137,14).Value = -10642.9998

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138,8).Value = 2017.2858
$ws.Cells.Item(138,9).Value = 1531.7
$ws.Cells.Item(138,11).Value = 4595.1
$ws.Cells.Item(138,13).Value = 544.8999999999996

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value = 3273.4285
$ws.Cells.Item(2,9).Value = 2028.8
$ws.Cells.Item(2,11).Value = 2028.8
$ws.Cells.Item(2,13).Value = -1915.8

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110,8).Value = 2621.1765
$ws.Cells.Item(110,9).Value = 1071.6666
$ws.Cells.Item(110,11).Value = 1071.6666
$ws.Cells.Item(110,13).Value = 973.3334

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116,8).Value = 3273.4285
$ws.Cells.Item(116,9).Value = 2028.8
$ws.Cells.Item(116,11).Value = 2028.8
$ws.Cells.Item(116,13).Value = 265.2

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value = 3273.4285
$ws.Cells.Item(3,9).Value = 2028.8
$ws.Cells.Item(3,11).Value = 2028.8
$ws.Cells.Item(3,13).Value = -1914.8

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134,8).Value = 5869.6943
$ws.Cells.Item(134,9).Value = 2760.6
$ws.Cells.Item(134,10).Value = 12935.818
$ws.Cells.Item(134,11).Value = 8281.799999999999
$ws.Cells.Item(134,12).Value = 38807.454
$ws.Cells.Item(134,13).Value = -5746.799999999999
$ws.Cells.Item(134,14).Value = -43877.454

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value = 4373.5835
$ws.Cells.Item(16,9).Value = 4269
$ws.Cells.Item(16,10).Value = 4520
$ws.Cells.Item(16,11).Value = 4269
$ws.Cells.Item(16,12).Value = 4520
$ws.Cells.Item(16,13).Value = -3982
$ws.Cells.Item(16,14).Value = -5094

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 5378462.5
$ws.Cells.Item(31,9).Value = 1413.641
$ws.Cells.Item(31,10).Value = 14496067
$ws.Cells.Item(31,11).Value = 1413.641
$ws.Cells.Item(31,12).Value = 14496067
$ws.Cells.Item(31,13).Value = -1118.641
$ws.Cells.Item(31,14).Value = -14496657

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34,8).Value = 5378462.5
$ws.Cells.Item(34,9).Value = 1413.641
$ws.Cells.Item(34,10).Value = 14496067
$ws.Cells.Item(34,11).Value = 1413.641
$ws.Cells.Item(34,12).Value = 14496067
$ws.Cells.Item(34,13).Value = -1211.641
$ws.Cells.Item(34,14).Value = -14496471

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94,8).Value = 2795
$ws.Cells.Item(94,9).Value = 3285.3333
$ws.Cells.Item(94,10).Value = 2304.6667
$ws.Cells.Item(94,11).Value = 3285.3333
$ws.Cells.Item(94,12).Value = 2304.6667
$ws.Cells.Item(94,13).Value = -2834.3333
$ws.Cells.Item(94,14).Value = -3206.6667

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105,8).Value = 2481.25
$ws.Cells.Item(105,9).Value = 2308.3333
$ws.Cells.Item(105,11).Value = 2308.3333
$ws.Cells.Item(105,13).Value = -561.3332999999998

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113,8).Value = 4373.5835
$ws.Cells.Item(113,9).Value = 4269
$ws.Cells.Item(113,10).Value = 4520
$ws.Cells.Item(113,11).Value = 4269
$ws.Cells.Item(113,12).Value = 4520
$ws.Cells.Item(113,13).Value = -2099
$ws.Cells.Item(113,14).Value = -8860

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80,8).Value = 3431
$ws.Cells.Item(80,9).Value = 0
$ws.Cells.Item(80,10).Value = 3431
$ws.Cells.Item(80,11).Value = 0
$ws.Cells.Item(80,12).Value = 10293
$ws.Cells.Item(80,13).ClearContents()
$ws.Cells.Item(80,14).Value = -12165

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83,8).Value = 3431
$ws.Cells.Item(83,9).Value = 0
$ws.Cells.Item(83,10).Value = 3431
$ws.Cells.Item(83,11).Value = 0
$ws.Cells.Item(83,12).Value = 30879
$ws.Cells.Item(83,13).ClearContents()
$ws.Cells.Item(83,14).Value = -40239

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113,8).Value = 543.32556
$ws.Cells.Item(113,9).Value = 443.86667
$ws.Cells.Item(113,10).Value = 596.6070999999999
$ws.Cells.Item(113,11).Value = 1331.60001
$ws.Cells.Item(113,12).Value = 1789.8213
$ws.Cells.Item(113,13).Value = 838.3999899999999
$ws.Cells.Item(113,14).Value = -6129.8213

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107,8).Value = 775.6667
$ws.Cells.Item(107,9).Value = 700.4
$ws.Cells.Item(107,10).Value = 829.4286
$ws.Cells.Item(107,11).Value = 700.4
$ws.Cells.Item(107,12).Value = 829.4286
$ws.Cells.Item(107,13).Value = 1219.6
$ws.Cells.Item(107,14).Value = -4669.4286

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 920.6667
$ws.Cells.Item(22,9).Value = 333.75
$ws.Cells.Item(22,10).Value = 1214.125
$ws.Cells.Item(22,11).Value = 333.75
$ws.Cells.Item(22,12).Value = 1214.125
$ws.Cells.Item(22,13).Value = -38.75
$ws.Cells.Item(22,14).Value = -1804.125

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27,8).Value = 920.6667
$ws.Cells.Item(27,9).Value = 333.75
$ws.Cells.Item(27,10).Value = 1214.125
$ws.Cells.Item(27,11).Value = 333.75
$ws.Cells.Item(27,12).Value = 1214.125
$ws.Cells.Item(27,13).Value = -226.75
$ws.Cells.Item(27,14).Value = -1428.125

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93,8).Value = 1759.64
$ws.Cells.Item(93,9).Value = 1604.85
$ws.Cells.Item(93,10).Value = 2378.8
$ws.Cells.Item(93,11).Value = 1604.85
$ws.Cells.Item(93,12).Value = 2378.8
$ws.Cells.Item(93,13).Value = -356.8499999999999
$ws.Cells.Item(93,14).Value = -4874.8

# WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41,8).Value = 21950.8
$ws.Cells.Item(41,9).Value = 0
$ws.Cells.Item(41,10).Value = 21950.8
$ws.Cells.Item(41,11).Value = 0
$ws.Cells.Item(41,12).Value = 21950.8
$ws.Cells.Item(41,13).ClearContents()
$ws.Cells.Item(41,14).Value = -22730.8
